# Update the weekly "Perejil" price-report sheet: every existing record
# shifts down by two rows (newest-first reordering as two new weekly
# records are folded in at the top), and two brand-new rows are appended
# at the bottom to keep the full history.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row target values for the variable columns:
#   D = Fecha (date serial), I = Calidad, J = Volumen,
#   K = Precio minimo, L = Precio maximo, M = Precio promedio ponderado,
#   P = Precio $/Kg
$rowsData = @{}
$rowsData[10] = @{ D=45069; I="Primera"; J=150; K=1200; L=1200; M=1200; P=1200 }
$rowsData[11] = @{ D=45069; I="Segunda"; J=150; K=1000; L=1000; M=1000; P=1000 }
$rowsData[12] = @{ D=45022; I="Primera"; J=230; K=1400; L=1500; M=1465; P=1465 }
$rowsData[13] = @{ D=45037; I="Primera"; J=100; K=1500; L=1500; M=1500; P=1500 }
$rowsData[14] = @{ D=45036; I="Primera"; J=250; K=1500; L=1500; M=1500; P=1500 }
$rowsData[15] = @{ D=44799; I="Primera"; J=160; K=750; L=850; M=800; P=800 }
$rowsData[16] = @{ D=44799; I="Segunda"; J=120; K=650; L=650; M=650; P=650 }
$rowsData[17] = @{ D=44804; I="Primera"; J=200; K=750; L=850; M=800; P=800 }
$rowsData[18] = @{ D=44804; I="Segunda"; J=200; K=650; L=650; M=650; P=650 }
$rowsData[19] = @{ D=44797; I="Primera"; J=240; K=750; L=850; M=800; P=800 }
$rowsData[20] = @{ D=44797; I="Segunda"; J=200; K=650; L=650; M=650; P=650 }
$rowsData[21] = @{ D=45044; I="Primera"; J=100; K=1500; L=1500; M=1500; P=1500 }
$rowsData[22] = @{ D=45054; I="Primera"; J=100; K=1500; L=1500; M=1500; P=1500 }
$rowsData[23] = @{ D=45054; I="Segunda"; J=100; K=1000; L=1000; M=1000; P=1000 }
$rowsData[24] = @{ D=45049; I="Primera"; J=100; K=1200; L=1200; M=1200; P=1200 }
$rowsData[25] = @{ D=45049; I="Segunda"; J=150; K=1000; L=1000; M=1000; P=1000 }
$rowsData[26] = @{ D=44761; I="Primera"; J=200; K=700; L=800; M=750; P=750 }
$rowsData[27] = @{ D=44761; I="Segunda"; J=150; K=600; L=600; M=600; P=600 }
$rowsData[28] = @{ D=44533; I="Primera"; J=100; K=2000; L=2200; M=2100; P=2100 }
$rowsData[29] = @{ D=45055; I="Primera"; J=150; K=1300; L=1300; M=1300; P=1300 }
$rowsData[30] = @{ D=45055; I="Segunda"; J=150; K=1000; L=1000; M=1000; P=1000 }
$rowsData[31] = @{ D=44818; I="Primera"; J=300; K=800; L=900; M=850; P=850 }
$rowsData[32] = @{ D=45002; I="Primera"; J=100; K=1200; L=1200; M=1200; P=1200 }
$rowsData[33] = @{ D=45033; I="Primera"; J=300; K=1500; L=1500; M=1500; P=1500 }
$rowsData[34] = @{ D=45068; I="Primera"; J=100; K=1200; L=1200; M=1200; P=1200 }
$rowsData[35] = @{ D=45068; I="Segunda"; J=100; K=1000; L=1000; M=1000; P=1000 }
$rowsData[36] = @{ D=44764; I="Primera"; J=200; K=700; L=800; M=750; P=750 }
$rowsData[37] = @{ D=44764; I="Segunda"; J=150; K=600; L=600; M=600; P=600 }
$rowsData[38] = @{ D=44754; I="Primera"; J=200; K=700; L=750; M=725; P=725 }
$rowsData[39] = @{ D=44882; I="Primera"; J=400; K=700; L=800; M=750; P=750 }
$rowsData[40] = @{ D=44882; I="Segunda"; J=300; K=600; L=600; M=600; P=600 }
$rowsData[41] = @{ D=45021; I="Primera"; J=200; K=1500; L=1500; M=1500; P=1500 }
$rowsData[42] = @{ D=45035; I="Primera"; J=150; K=1500; L=1500; M=1500; P=1500 }
$rowsData[43] = @{ D=44610; I="Primera"; J=100; K=600; L=650; M=625; P=625 }
$rowsData[44] = @{ D=44837; I="Primera"; J=200; K=700; L=800; M=750; P=750 }
$rowsData[45] = @{ D=44837; I="Segunda"; J=150; K=600; L=600; M=600; P=600 }
$rowsData[46] = @{ D=44791; I="Primera"; J=240; K=750; L=800; M=775; P=775 }
$rowsData[47] = @{ D=44791; I="Segunda"; J=250; K=650; L=650; M=650; P=650 }
$rowsData[48] = @{ D=45030; I="Primera"; J=300; K=1500; L=1500; M=1500; P=1500 }
$rowsData[49] = @{ D=44831; I="Primera"; J=300; K=700; L=800; M=750; P=750 }
$rowsData[50] = @{ D=44831; I="Segunda"; J=200; K=600; L=600; M=600; P=600 }
$rowsData[51] = @{ D=45041; I="Primera"; J=200; K=1500; L=1500; M=1500; P=1500 }
$rowsData[52] = @{ D=44811; I="Primera"; J=300; K=750; L=850; M=800; P=800 }
$rowsData[53] = @{ D=45016; I="Primera"; J=100; K=1500; L=1500; M=1500; P=1500 }
$rowsData[54] = @{ D=45043; I="Primera"; J=300; K=1500; L=1500; M=1500; P=1500 }
$rowsData[55] = @{ D=45027; I="Primera"; J=200; K=1500; L=1500; M=1500; P=1500 }
$rowsData[56] = @{ D=45001; I="Primera"; J=150; K=1300; L=1300; M=1300; P=1300 }
$rowsData[57] = @{ D=44532; I="Primera"; J=60; K=2000; L=2200; M=2100; P=2100 }
$rowsData[58] = @{ D=45020; I="Primera"; J=200; K=1400; L=1500; M=1475; P=1475 }
$rowsData[59] = @{ D=44859; I="Primera"; J=300; K=700; L=800; M=750; P=750 }
$rowsData[60] = @{ D=44859; I="Segunda"; J=200; K=600; L=600; M=600; P=600 }
$rowsData[61] = @{ D=44624; I="Primera"; J=120; K=650; L=700; M=675; P=675 }
$rowsData[62] = @{ D=45014; I="Primera"; J=60; K=1500; L=1500; M=1500; P=1500 }
$rowsData[63] = @{ D=44999; I="Primera"; J=200; K=1500; L=1500; M=1500; P=1500 }

for ($r = 10; $r -le 63; $r++) {
    $row = $rowsData[$r]
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 16).Value = $row.P
}

# Rows 62 and 63 are brand new records appended at the bottom of the
# table; fill in the columns that are constant across every row in this
# sheet (they were blank before this edit), and give the new "Fecha"
# cells the same date number format used by the rest of column D.
$ws.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(62, 1).Value = 7
$ws.Cells.Item(62, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(62, 3).Value = "Ñuble"
$ws.Cells.Item(62, 5).Value = 16
$ws.Cells.Item(62, 6).Value = 100112044
$ws.Cells.Item(62, 7).Value = "Perejil"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(62, 15).Value = "Región del Maule"
$ws.Cells.Item(62, 17).Value = 1
$ws.Cells.Item(62, 18).Value = "Hortaliza"

$ws.Cells.Item(63, 1).Value = 7
$ws.Cells.Item(63, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(63, 3).Value = "Ñuble"
$ws.Cells.Item(63, 5).Value = 16
$ws.Cells.Item(63, 6).Value = 100112044
$ws.Cells.Item(63, 7).Value = "Perejil"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(63, 15).Value = "Región del Maule"
$ws.Cells.Item(63, 17).Value = 1
$ws.Cells.Item(63, 18).Value = "Hortaliza"

Write-Output "Perejil sheet updated: rows 10-63 rewritten, dimension now A1:R63."
